# The source data added one more weekly observation for "Apio" (Vega Central
# Mapocho de Santiago). This shifts every existing data row from 510 onward
# down by one, and the freshly opened row 510 receives the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 510..594 down to 511..595, opening up a blank row 510.
$ws.Rows(510).Insert()

# Populate the newly inserted row 510 with the new observation.
$ws.Cells.Item(510, 1).Value  = 9
$ws.Cells.Item(510, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(510, 3).Value  = "Metropolitana"
$ws.Cells.Item(510, 4).Value  = 45258
$ws.Cells.Item(510, 5).Value  = 13
$ws.Cells.Item(510, 6).Value  = 100112017
$ws.Cells.Item(510, 7).Value  = "Apio"
$ws.Cells.Item(510, 8).Value  = "Americana (o)"
$ws.Cells.Item(510, 9).Value  = "Segunda"
$ws.Cells.Item(510, 10).Value = 60
$ws.Cells.Item(510, 11).Value = 10000
$ws.Cells.Item(510, 12).Value = 10000
$ws.Cells.Item(510, 13).Value = 10000
$ws.Cells.Item(510, 14).Value = "$/docena de matas"
$ws.Cells.Item(510, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(510, 16).Value = 1667
$ws.Cells.Item(510, 17).Value = 6
$ws.Cells.Item(510, 18).Value = "Hortaliza"
